$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy border+center-alignment style from column A into columns B/C
#     for every data row (2..31) so new cells inherit the bordered "s=2" look. ---
for ($r = 2; $r -le 31; $r++) {
    $ws.Range("A$r").Copy()
    $ws.Range("B$r").PasteSpecial(-4122)
    $ws.Range("C$r").PasteSpecial(-4122)
}

# --- Step 2: populate cell values in the precise order needed so new shared
#     strings land at the same indices as the target workbook. ---

# 1) "yes"
$ws.Range("B3").Value = "yes"

# 2) whole UOM column (C2..C31, skipping the blank C3) top-to-bottom
$ws.Range("C2").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("C5").Value = "DD/MM/YYYY"
$ws.Range("C6").Value = "-"
$ws.Range("C7").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("C9").Value = "dS/m"
$ws.Range("C10").Value = "-"
$ws.Range("C11").Value = "ppm"
$ws.Range("C12").Value = "ppm"
$ws.Range("C13").Value = "ppm"
$ws.Range("C14").Value = "ppm"
$ws.Range("C15").Value = "ppm"
$ws.Range("C16").Value = "ppm"
$ws.Range("C17").Value = "ppm"
$ws.Range("C18").Value = "ppm"
$ws.Range("C19").Value = "ppm"
$ws.Range("C20").Value = "ppm"
$ws.Range("C21").Value = "ppm"
$ws.Range("C22").Value = "ppm"
$ws.Range("C23").Value = "ppm"
$ws.Range("C24").Value = "ppm"
$ws.Range("C25").Value = "ppm"
$ws.Range("C26").Value = "ppm"
$ws.Range("C27").Value = "ppm"
$ws.Range("C28").Value = "ppm"
$ws.Range("C29").Value = "ppm"
$ws.Range("C30").Value = "ppm"
$ws.Range("C31").Value = "CFU/ml"

# 3) "Irrigation", then "previous"
$ws.Range("B7").Value = "Irrigation"
$ws.Range("B6").Value = "previous"

# --- Step 3: fill in the remaining (numeric) Value-column cells ---
$ws.Range("B2").Value = 100
$ws.Range("B4").Value = 15
$ws.Range("B5").Value = 43235
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 0.3
$ws.Range("B10").Value = 7.7
$ws.Range("B11").Value = 1.8
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 5
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 3
$ws.Range("B16").Value = 48
$ws.Range("B17").Value = 8.2
$ws.Range("B18").Value = 28
$ws.Range("B19").Value = 0
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 0
$ws.Range("B22").Value = 0
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
$ws.Range("B25").Value = 0
$ws.Range("B26").Value = 0
$ws.Range("B27").Value = 0
$ws.Range("B28").Value = 0
$ws.Range("B29").Value = 0
$ws.Range("B30").Value = 0
$ws.Range("B31").Value = 0

# --- Step 4: date format for the Sample_Date value cell ---
$ws.Range("B5").NumberFormat = "mm-dd-yy"

# --- Step 5: data validation (dropdown list) on the Source value cell ---
$ws.Range("B7").Validation.Add(3, 1, 1, "=#REF!")

$v = $ws.Range("B2").Value()
Write-Output "Edit complete. B2=$v"
